$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Update the date heading
$d.Content.Find.Execute("2024-06-26 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-27 Thursday", 2)

# Row 1
$t.Cell(1,1).Range.Text = "29÷8="
$t.Cell(1,2).Range.Text = "48÷5="
$t.Cell(1,3).Range.Text = "64÷2="
$t.Cell(1,4).Range.Text = "13÷2="
$t.Cell(1,5).Range.Text = "21÷7="

# Row 5
$t.Cell(5,1).Range.Text = "11÷6="
$t.Cell(5,2).Range.Text = "92÷2="
$t.Cell(5,3).Range.Text = "64÷2="
$t.Cell(5,4).Range.Text = "15÷9="
$t.Cell(5,5).Range.Text = "23÷6="

# Row 9
$t.Cell(9,1).Range.Text = "84÷8="
$t.Cell(9,2).Range.Text = "79÷9="
$t.Cell(9,3).Range.Text = "10÷4="
$t.Cell(9,4).Range.Text = "37÷5="
$t.Cell(9,5).Range.Text = "43÷9="

# Row 13
$t.Cell(13,1).Range.Text = "41÷5="
$t.Cell(13,2).Range.Text = "16÷9="
$t.Cell(13,3).Range.Text = "58÷4="
$t.Cell(13,4).Range.Text = "66÷6="
$t.Cell(13,5).Range.Text = "65÷3="

# Row 17
$t.Cell(17,1).Range.Text = "91÷2="
$t.Cell(17,2).Range.Text = "51÷3="
$t.Cell(17,3).Range.Text = "31÷3="
$t.Cell(17,4).Range.Text = "89÷7="
$t.Cell(17,5).Range.Text = "75÷9="
